# Adds a "Price" column (C) to the Products sheet, with currency-formatted
# values for each product row, and resizes the columns to make room for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Give column C (header + data rows) the same look (font/border/
#        alignment) as the existing A/B columns before we put data in it,
#        by copying the formatting from the header cell A1. ---
$ws.Range("A1:A7").Copy()
$ws.Range("C1:C7").PasteSpecial(-4122)

# --- 2. Header ---
$ws.Range("C1").Value = "Price"

# --- 3. Data values (Price for each of the 6 products) ---
$ws.Range("C2").Value = 29.99
$ws.Range("C3").Value = 9.99
$ws.Range("C4").Value = 15.99
$ws.Range("C5").Value = 49.99
$ws.Range("C6").Value = 7.99
$ws.Range("C7").Value = 15.99

# --- 4. Number format for the price cells: built-in currency format (id 7) ---
$ws.Range("C2:C7").NumberFormat = "$#,##0.00_);($#,##0.00)"

# --- 5. Column widths: keep column A the same, narrow column B slightly and
#        size the new column C, to fit the new layout (values chosen so the
#        saved/stored width matches the target layout as closely as the
#        ColumnWidth property allows). ---
$ws.Columns.Item(1).ColumnWidth = 28.5
$ws.Columns.Item(2).ColumnWidth = 150.83333333333334
$ws.Columns.Item(3).ColumnWidth = 47.666666666666664
